# Scheduled-runner update: refresh Leve profit figures (currentAveragePrice /
# NQ / HQ and derived Leve profit columns) across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 129.625
$ws.Range("I12").Value = 159.2
$ws.Range("J12").Value = 80.333336
$ws.Range("K12").Value = 159.2
$ws.Range("L12").Value = 80.333336
$ws.Range("M12").Value = 10.80000000000001
$ws.Range("N12").Value = -420.333336

$ws.Range("H17").Value = 3285.25
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3285.25
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 9855.75
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -10191.75

$ws.Range("H40").Value = 5500
$ws.Range("J40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("N40").Value = -1350

$ws.Range("H64").Value = 6499.75
$ws.Range("I64").Value = 6499.75
$ws.Range("K64").Value = 6499.75
$ws.Range("M64").Value = -6251.75

$ws.Range("H67").Value = 6499.75
$ws.Range("I67").Value = 6499.75
$ws.Range("K67").Value = 6499.75
$ws.Range("M67").Value = -5641.75

$ws.Range("H92").Value = 911.06665
$ws.Range("I92").Value = 882.1539
$ws.Range("K92").Value = 882.1539
$ws.Range("M92").Value = 365.8461

$ws.Range("H132").Value = 2431.7
$ws.Range("I132").Value = 2431.7
$ws.Range("K132").Value = 7295.099999999999
$ws.Range("M132").Value = -4765.099999999999

$ws.Range("H138").Value = 8315.799999999999
$ws.Range("J138").Value = 9644.75
$ws.Range("L138").Value = 28934.25
$ws.Range("N138").Value = -39214.25

$ws.Range("H141").Value = 2674.75
$ws.Range("I141").Value = 899.6667
$ws.Range("K141").Value = 2699.0001
$ws.Range("M141").Value = 2480.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2381.2
$ws.Range("I2").Value = 3202
$ws.Range("J2").Value = 1150
$ws.Range("K2").Value = 3202
$ws.Range("L2").Value = 1150
$ws.Range("M2").Value = -3089
$ws.Range("N2").Value = -1376

$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("M19").Value = -271

$ws.Range("H74").Value = 9611.788
$ws.Range("I74").Value = 9757.204
$ws.Range("J74").Value = 8812
$ws.Range("K74").Value = 9757.204
$ws.Range("L74").Value = 8812
$ws.Range("M74").Value = -8883.204
$ws.Range("N74").Value = -10560

$ws.Range("H77").Value = 9611.788
$ws.Range("I77").Value = 9757.204
$ws.Range("J77").Value = 8812
$ws.Range("K77").Value = 48786.02
$ws.Range("L77").Value = 44060
$ws.Range("M77").Value = -44418.02
$ws.Range("N77").Value = -52796

$ws.Range("H116").Value = 2381.2
$ws.Range("I116").Value = 3202
$ws.Range("J116").Value = 1150
$ws.Range("K116").Value = 3202
$ws.Range("L116").Value = 1150
$ws.Range("M116").Value = -908
$ws.Range("N116").Value = -5738

$ws.Range("H132").Value = 1811.625
$ws.Range("I132").Value = 1854.2727
$ws.Range("K132").Value = 5562.8181
$ws.Range("M132").Value = -3032.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2381.2
$ws.Range("I3").Value = 3202
$ws.Range("J3").Value = 1150
$ws.Range("K3").Value = 3202
$ws.Range("L3").Value = 1150
$ws.Range("M3").Value = -3088
$ws.Range("N3").Value = -1378

$ws.Range("H94").Value = 3902.25
$ws.Range("J94").Value = 800
$ws.Range("L94").Value = 800
$ws.Range("N94").Value = -1702

$ws.Range("H99").Value = 2467.6
$ws.Range("I99").Value = 2772
$ws.Range("J99").Value = 2011
$ws.Range("K99").Value = 2772
$ws.Range("L99").Value = 2011
$ws.Range("M99").Value = -1274
$ws.Range("N99").Value = -5007

$ws.Range("H107").Value = 656.0909
$ws.Range("I107").Value = 668.55554
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 668.55554
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 1251.44446
$ws.Range("N107").Value = -4440

$ws.Range("H134").Value = 3540.4
$ws.Range("I134").Value = 3563
$ws.Range("K134").Value = 10689
$ws.Range("M134").Value = -8154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8333.333000000001
$ws.Range("I99").Value = 8333.333000000001
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8333.333000000001
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -6835.333000000001
$ws.Range("N99").ClearContents()

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H126").Value = 8333.333000000001
$ws.Range("I126").Value = 8333.333000000001
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 24999.999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -22529.999
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3594.4167
$ws.Range("I132").Value = 3353.6
$ws.Range("J132").Value = 4798.5
$ws.Range("K132").Value = 10060.8
$ws.Range("L132").Value = 14395.5
$ws.Range("M132").Value = -7530.799999999999
$ws.Range("N132").Value = -19455.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4999.25
$ws.Range("J81").Value = 5499
$ws.Range("L81").Value = 16497
$ws.Range("N81").Value = -18743

$ws.Range("H84").Value = 4999.25
$ws.Range("J84").Value = 5499
$ws.Range("L84").Value = 49491
$ws.Range("N84").Value = -60723

$ws.Range("H86").Value = 851.5454999999999
$ws.Range("I86").Value = 771
$ws.Range("J86").Value = 881.75
$ws.Range("K86").Value = 2313
$ws.Range("L86").Value = 2645.25
$ws.Range("M86").Value = -1127
$ws.Range("N86").Value = -5017.25

$ws.Range("H89").Value = 851.5454999999999
$ws.Range("I89").Value = 771
$ws.Range("J89").Value = 881.75
$ws.Range("K89").Value = 6939
$ws.Range("L89").Value = 7935.75
$ws.Range("M89").Value = -1011
$ws.Range("N89").Value = -19791.75

$ws.Range("H122").Value = 3601.739
$ws.Range("I122").Value = 1516.6666
$ws.Range("J122").Value = 3914.5
$ws.Range("K122").Value = 13649.9994
$ws.Range("L122").Value = 35230.5
$ws.Range("M122").Value = -11199.9994
$ws.Range("N122").Value = -40130.5

$ws.Range("H131").Value = 2498.8333
$ws.Range("I131").Value = 2995
$ws.Range("J131").Value = 2436.8125
$ws.Range("K131").Value = 8985
$ws.Range("L131").Value = 7310.4375
$ws.Range("M131").Value = -3945
$ws.Range("N131").Value = -17390.4375

$ws.Range("H141").Value = 2549.5
$ws.Range("I141").Value = 2549.5
$ws.Range("K141").Value = 7648.5
$ws.Range("M141").Value = -2468.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2357.8
$ws.Range("I102").Value = 2357.8
$ws.Range("K102").Value = 2357.8
$ws.Range("M102").Value = -735.8000000000002

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()

$ws.Range("H106").Value = 51132.332
$ws.Range("J106").Value = 51132.332
$ws.Range("L106").Value = 51132.332
$ws.Range("N106").Value = -53656.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1196.6364
$ws.Range("I100").Value = 826.875
$ws.Range("K100").Value = 1653.75
$ws.Range("M100").Value = -1112.75

$ws.Range("H132").Value = 868.5
$ws.Range("I132").Value = 902.2
$ws.Range("J132").Value = 700
$ws.Range("K132").Value = 2706.6
$ws.Range("L132").Value = 2100
$ws.Range("M132").Value = -176.6000000000004
$ws.Range("N132").Value = -7160
